$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3007
$ws.Range("K3").Value = 2961
$ws.Range("E4").Value = 758
$ws.Range("K4").Value = 606
$ws.Range("K5").Value = 192
$ws.Range("K6").Value = 3549
$ws.Range("E7").Value = 9484
$ws.Range("K7").Value = 10315

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 303
$ws.Range("K8").Value = 681
$ws.Range("K10").Value = 58
$ws.Range("K11").Value = 216
$ws.Range("K15").Value = 104
$ws.Range("K16").Value = 33
$ws.Range("K18").Value = 70
$ws.Range("K19").Value = 313
$ws.Range("K20").Value = 239
$ws.Range("K23").Value = 92
$ws.Range("K24").Value = 35
$ws.Range("K27").Value = 107
$ws.Range("K29").Value = 539
$ws.Range("K31").Value = 116
$ws.Range("K33").Value = 410
$ws.Range("K36").Value = 121
$ws.Range("K37").Value = 341
$ws.Range("E41").Value = 57
$ws.Range("K41").Value = 88
$ws.Range("K42").Value = 356
$ws.Range("K43").Value = 92
$ws.Range("K44").Value = 101
$ws.Range("K48").Value = 120
$ws.Range("K54").Value = 199
$ws.Range("K55").Value = 108
$ws.Range("K57").Value = 28
$ws.Range("K58").Value = 2
$ws.Range("K60").Value = 64
$ws.Range("K63").Value = 38
$ws.Range("K66").Value = 37
$ws.Range("K67").Value = 409
$ws.Range("K75").Value = 38
$ws.Range("K76").Value = 161
$ws.Range("K77").Value = 72
$ws.Range("K78").Value = 136
$ws.Range("K79").Value = 265
$ws.Range("K83").Value = 225
$ws.Range("K84").Value = 73
$ws.Range("K85").Value = 493
$ws.Range("K88").Value = 115
$ws.Range("K91").Value = 105
$ws.Range("K94").Value = 122
$ws.Range("K95").Value = 169
$ws.Range("K99").Value = 185
$ws.Range("E101").Value = 9484
$ws.Range("K101").Value = 10315

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 92
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 303

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 179
$ws.Range("K3").Value = 169
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 493

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 196
$ws.Range("K3").Value = 206
$ws.Range("K4").Value = 37
$ws.Range("K6").Value = 226
$ws.Range("K7").Value = 681

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 148
$ws.Range("K6").Value = 118
$ws.Range("K7").Value = 410

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 169

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 89
$ws.Range("K3").Value = 117
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 341

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 67
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 185

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 134
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 409

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 183
$ws.Range("K6").Value = 165
$ws.Range("K7").Value = 539

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 106
$ws.Range("K7").Value = 313

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("E4").Value = 5
$ws.Range("K4").Value = 5
$ws.Range("E7").Value = 57
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 91
$ws.Range("K6").Value = 136
$ws.Range("K7").Value = 356

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 26
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 89
$ws.Range("K3").Value = 94
$ws.Range("K5").Value = 9
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 83
$ws.Range("K3").Value = 67
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K5").Value = 20
$ws.Range("K6").Value = 37

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 27
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 2
